$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row: Right count 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row: Right total 54 -> 90
$ws.Range("B12").Value = 90

# Update the correct/total marks text 52/84 -> 90/140
$ws.Range("E12").Value = "90/140"
